$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "growmat_easy_hydro"

# --- Header row ---
$ws2.Range("B1").Value = 'item'
$ws2.Range("C1").Value = 'pcs'
$ws2.Range("D1").Value = 'price per pcs'
$ws2.Range("E1").Value = 'price total'
$ws2.Range("F1").Value = 'coeficient'
$ws2.Range("G1").Value = 'price'
$ws2.Range("H1").Value = 'sales price'
$ws2.Range("I1").Value = 'source'
$ws2.Range("J1").Value = 'vat'
$ws2.Range("L1").Value = 'descr'
$ws2.Range("M1").Value = 'link'

# --- Data rows 3-30 ---
# row 3: box
$ws2.Range("B3").Value = 'box'
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 280
$ws2.Range("E3").Formula = "=D3*C3"
$ws2.Range("F3").Value = 1
$ws2.Range("G3").Formula = "=F3*E3"
$ws2.Range("I3").Value = 'cz'
$ws2.Range("J3").Value = 1
$ws2.Range("M3").Value = 'hornbach'

# row 4: 433MHz socket
$ws2.Range("B4").Value = '433MHz socket'
$ws2.Range("C4").Value = 1
$ws2.Range("D4").Value = 700
$ws2.Range("E4").Formula = "=D4*C4"
$ws2.Range("F4").Value = 1
$ws2.Range("G4").Formula = "=F4*E4"
$ws2.Range("I4").Value = 'cz'
$ws2.Range("J4").Value = 1
$ws2.Range("M4").Value = 'http://www.giganto.cz/elektronika/dalkove-ovladane-bezdratove-zasuvky-3600w'

# row 5: 433MHz tx
$ws2.Range("B5").Value = '433MHz tx'
$ws2.Range("C5").Value = 1
$ws2.Range("D5").Value = 40
$ws2.Range("E5").Formula = "=D5*C5"
$ws2.Range("F5").Value = 1
$ws2.Range("G5").Formula = "=F5*E5"
$ws2.Range("I5").Value = 'cz'
$ws2.Range("J5").Value = 0
$ws2.Range("M5").Value = 'https://www.postavrobota.cz/Vysilac-433MHz-ASK-antena-d663.htm'

# row 6: display
$ws2.Range("B6").Value = 'display'
$ws2.Range("C6").Value = 1
$ws2.Range("D6").Value = 160
$ws2.Range("E6").Formula = "=D6*C6"
$ws2.Range("F6").Value = 1
$ws2.Range("G6").Formula = "=F6*E6"
$ws2.Range("I6").Value = 'cz'
$ws2.Range("J6").Value = 0
$ws2.Range("M6").Value = 'https://www.postavrobota.cz/I2C-LCD-displej-znakovy-16x2-modry-d333.htm'

# row 7: keyboard
$ws2.Range("B7").Value = 'keyboard'
$ws2.Range("C7").Value = 1
$ws2.Range("D7").Value = 40
$ws2.Range("E7").Formula = "=D7*C7"
$ws2.Range("F7").Value = 1
$ws2.Range("G7").Formula = "=F7*E7"
$ws2.Range("I7").Value = 'cz'
$ws2.Range("J7").Value = 0
$ws2.Range("M7").Value = 'https://www.postavrobota.cz/Membranova-klavesnice-4x3-samolepici-d137.htm'

# row 8: i2c expander
$ws2.Range("B8").Value = 'i2c expander'
$ws2.Range("C8").Value = 1
$ws2.Range("D8").Value = 50
$ws2.Range("E8").Formula = "=D8*C8"
$ws2.Range("F8").Value = 2
$ws2.Range("G8").Formula = "=F8*E8"
$ws2.Range("I8").Value = 'ebay'
$ws2.Range("M8").Value = 'http://www.ebay.com/itm/PCF8574-PCF8574T-I2C-8-Bit-IO-GPIO-Expander-Module-for-Arduino-Raspberry-Pi-UK-/272432637606?var=&hash=item3f6e4026a6:m:mQFPRLxSj-zVxa61Qe6YVhQ'

# row 9: arduino mega
$ws2.Range("B9").Value = 'arduino mega'
$ws2.Range("C9").Value = 1
$ws2.Range("D9").Value = 400
$ws2.Range("E9").Formula = "=D9*C9"
$ws2.Range("F9").Value = 1
$ws2.Range("G9").Formula = "=F9*E9"
$ws2.Range("I9").Value = 'cz'
$ws2.Range("J9").Value = 0
$ws2.Range("M9").Value = 'https://www.postavrobota.cz/Dccduino-Mega-ATmega2560-Arduino-kompatibilni-d100.htm'

# row 10: battery
$ws2.Range("B10").Value = 'battery'
$ws2.Range("C10").Value = 1
$ws2.Range("D10").Value = 50
$ws2.Range("E10").Formula = "=D10*C10"
$ws2.Range("F10").Value = 2
$ws2.Range("G10").Formula = "=F10*E10"
$ws2.Range("I10").Value = 'ebay'
$ws2.Range("M10").Value = 'http://www.ebay.com/itm/4pcs-3-7V-18650-9900mah-Li-ion-Rechargeable-Battery-For-LED-Flashlight-Torch-LO-/161933546105?hash=item25b3fdd279:g:JuIAAOSwZG9WhJCM'

# row 11: charger
$ws2.Range("B11").Value = 'charger'
$ws2.Range("C11").Value = 1
$ws2.Range("D11").Value = 45
$ws2.Range("E11").Formula = "=D11*C11"
$ws2.Range("F11").Value = 1
$ws2.Range("G11").Formula = "=F11*E11"
$ws2.Range("I11").Value = 'cz'
$ws2.Range("J11").Value = 0
$ws2.Range("M11").Value = 'https://www.postavrobota.cz/Mini-nabijecka-Li-ion-Li-po-clanku-s-ochranou-baterie-d10.htm'

# row 12: dc step-up
$ws2.Range("B12").Value = 'dc step-up'
$ws2.Range("C12").Value = 1
$ws2.Range("D12").Value = 35
$ws2.Range("E12").Formula = "=D12*C12"
$ws2.Range("F12").Value = 1
$ws2.Range("G12").Formula = "=F12*E12"
$ws2.Range("I12").Value = 'cz'
$ws2.Range("J12").Value = 0
$ws2.Range("M12").Value = 'https://www.postavrobota.cz/Mini-nastavitelny-zdroj-zvysujici-napeti-2A-d136.htm'

# row 13: gsm
$ws2.Range("B13").Value = 'gsm'
$ws2.Range("C13").Value = 1
$ws2.Range("D13").Value = 300
$ws2.Range("E13").Formula = "=D13*C13"
$ws2.Range("F13").Value = 2
$ws2.Range("G13").Formula = "=F13*E13"
$ws2.Range("I13").Value = 'ebay'
$ws2.Range("M13").Value = 'http://www.ebay.com/itm/SIM800L-GPRS-GSM-SIM-Board-Quadband-QUAD-BAND-L-shape-Antenna-for-Arduino-/281958541036?hash=item41a609e6ec:g:5JUAAOSwMVdYFv3B'

# row 14: ph probe
$ws2.Range("B14").Value = 'ph probe'
$ws2.Range("C14").Value = 1
$ws2.Range("D14").Value = 700
$ws2.Range("E14").Formula = "=D14*C14"
$ws2.Range("F14").Value = 2
$ws2.Range("G14").Formula = "=F14*E14"
$ws2.Range("I14").Value = 'ebay'
$ws2.Range("M14").Value = 'http://www.ebay.com/itm/Liquid-PH0-14-Value-Detect-Sensor-Module-PH-Electrode-Probe-BNC-for-Arduino-G-/262617141991?hash=item3d253392e7:g:T0IAAOSwxg5X0kt3'

# row 15: ph modul
$ws2.Range("B15").Value = 'ph modul'
$ws2.Range("C15").Value = 1
$ws2.Range("D15").Value = 0
$ws2.Range("E15").Formula = "=D15*C15"
$ws2.Range("F15").Value = 1
$ws2.Range("G15").Formula = "=F15*E15"

# row 16: ec probe
$ws2.Range("B16").Value = 'ec probe'
$ws2.Range("C16").Value = 1
$ws2.Range("D16").Value = 700
$ws2.Range("E16").Formula = "=D16*C16"
$ws2.Range("F16").Value = 2
$ws2.Range("G16").Formula = "=F16*E16"
$ws2.Range("I16").Value = 'ebay'
$ws2.Range("M16").Value = 'http://www.ebay.com/itm/E201WM-Conductivity-COND-EC-electrode-Conductivity-sensor-probe-BNC-connector-/141936907599?hash=item210c19554f:g:eRUAAOSwEK9Txo38'

# row 17: ec modul
$ws2.Range("B17").Value = 'ec modul'
$ws2.Range("C17").Value = 1
$ws2.Range("D17").Value = 100
$ws2.Range("E17").Formula = "=D17*C17"
$ws2.Range("F17").Value = 1
$ws2.Range("G17").Formula = "=F17*E17"
$ws2.Range("I17").Value = 'cz'
$ws2.Range("J17").Value = 1
$ws2.Range("M17").Value = 'gsm'

# row 18: ec dc-dc isolated
$ws2.Range("B18").Value = 'ec dc-dc isolated'
$ws2.Range("C18").Value = 1
$ws2.Range("D18").Value = 100
$ws2.Range("E18").Formula = "=D18*C18"
$ws2.Range("F18").Value = 1
$ws2.Range("G18").Formula = "=F18*E18"
$ws2.Range("I18").Value = 'cz'
$ws2.Range("J18").Value = 0
$ws2.Range("M18").Value = 'https://www.postavrobota.cz/DC-DC-5V-5V-izolovany-zdroj-1W-d279.htm'

# row 19: thermo/humidyty meter air
$ws2.Range("B19").Value = 'thermo/humidyty meter air'
$ws2.Range("C19").Value = 1
$ws2.Range("D19").Value = 150
$ws2.Range("E19").Formula = "=D19*C19"
$ws2.Range("F19").Value = 2
$ws2.Range("G19").Formula = "=F19*E19"
$ws2.Range("I19").Value = 'ebay'

# row 20: light sensor
$ws2.Range("B20").Value = 'light sensor'
$ws2.Range("C20").Value = 1
$ws2.Range("D20").Value = 50
$ws2.Range("E20").Formula = "=D20*C20"
$ws2.Range("F20").Value = 1
$ws2.Range("G20").Formula = "=F20*E20"
$ws2.Range("I20").Value = 'cz'
$ws2.Range("J20").Value = 1
$ws2.Range("M20").Value = 'gsm'

# row 21: thermometer water
$ws2.Range("B21").Value = 'thermometer water'
$ws2.Range("C21").Value = 1
$ws2.Range("D21").Value = 150
$ws2.Range("E21").Formula = "=D21*C21"
$ws2.Range("F21").Value = 1
$ws2.Range("G21").Formula = "=F21*E21"
$ws2.Range("I21").Value = 'cz'
$ws2.Range("J21").Value = 0
$ws2.Range("M21").Value = 'https://www.postavrobota.cz/Digitalni-teplotni-sonda-DS18B20-2m-d370.htm'

# row 22: ultrasonic distance meter
$ws2.Range("B22").Value = 'ultrasonic distance meter'
$ws2.Range("C22").Value = 1
$ws2.Range("D22").Value = 200
$ws2.Range("E22").Formula = "=D22*C22"
$ws2.Range("F22").Value = 2
$ws2.Range("G22").Formula = "=F22*E22"
$ws2.Range("I22").Value = 'ebay'
$ws2.Range("M22").Value = 'http://www.ebay.com/itm/Ultrasonic-Module-Distance-Measuring-Transducer-Sensor-Waterproof-Perfect-/272041782549?hash=item3f56f42d15:g:mMYAAOSw5ZBWQbkT'

# row 23: level switch
$ws2.Range("B23").Value = 'level switch'
$ws2.Range("C23").Value = 1
$ws2.Range("D23").Value = 100
$ws2.Range("E23").Formula = "=D23*C23"
$ws2.Range("F23").Value = 1
$ws2.Range("G23").Formula = "=F23*E23"

# row 24: usb cable
$ws2.Range("B24").Value = 'usb cable'
$ws2.Range("C24").Value = 1
$ws2.Range("D24").Value = 50
$ws2.Range("E24").Formula = "=D24*C24"
$ws2.Range("F24").Value = 1
$ws2.Range("G24").Formula = "=F24*E24"

# row 25: power source
$ws2.Range("B25").Value = 'power source'
$ws2.Range("C25").Value = 1
$ws2.Range("D25").Value = 100
$ws2.Range("E25").Formula = "=D25*C25"
$ws2.Range("F25").Value = 1
$ws2.Range("G25").Formula = "=F25*E25"

# row 26: rtc
$ws2.Range("B26").Value = 'rtc'
$ws2.Range("C26").Value = 1
$ws2.Range("D26").Value = 50
$ws2.Range("E26").Formula = "=D26*C26"
$ws2.Range("F26").Value = 1
$ws2.Range("G26").Formula = "=F26*E26"
$ws2.Range("M26").Value = 'https://www.postavrobota.cz/RTC-modul-realneho-casu-DS1307-baterie-d123.htm'

# row 27: pcb ps
$ws2.Range("B27").Value = 'pcb ps'
$ws2.Range("C27").Value = 1
$ws2.Range("E27").Formula = "=D27*C27"
$ws2.Range("F27").Value = 1
$ws2.Range("G27").Formula = "=F27*E27"

# row 28: pcb ec
$ws2.Range("B28").Value = 'pcb ec'
$ws2.Range("C28").Value = 1
$ws2.Range("E28").Formula = "=D28*C28"
$ws2.Range("F28").Value = 1
$ws2.Range("G28").Formula = "=F28*E28"

# row 29: pcb ?
$ws2.Range("B29").Value = 'pcb ?'
$ws2.Range("C29").Value = 1
$ws2.Range("E29").Formula = "=D29*C29"
$ws2.Range("F29").Value = 1
$ws2.Range("G29").Formula = "=F29*E29"

# row 30: esp8266
$ws2.Range("B30").Value = 'esp8266'
$ws2.Range("C30").Value = 1
$ws2.Range("D30").Value = 140
$ws2.Range("E30").Formula = "=D30*C30"
$ws2.Range("F30").Value = 1
$ws2.Range("G30").Formula = "=F30*E30"
$ws2.Range("I30").Value = 'cz'
$ws2.Range("J30").Value = 0
$ws2.Range("M30").Value = 'https://www.postavrobota.cz/WiFi-RS232-AT-modul-ESP8266-2-4GHz-SoC-d159.htm'

# --- Rows 31-45 (blank items, just counters/formulas) ---
$ws2.Range("C31").Value = 1
$ws2.Range("E31").Formula = "=D31*C31"
$ws2.Range("F31").Value = 1
$ws2.Range("G31").Formula = "=F31*E31"
$ws2.Range("C32").Value = 1
$ws2.Range("E32").Formula = "=D32*C32"
$ws2.Range("F32").Value = 1
$ws2.Range("G32").Formula = "=F32*E32"
$ws2.Range("C33").Value = 1
$ws2.Range("E33").Formula = "=D33*C33"
$ws2.Range("F33").Value = 1
$ws2.Range("G33").Formula = "=F33*E33"
$ws2.Range("C34").Value = 1
$ws2.Range("E34").Formula = "=D34*C34"
$ws2.Range("F34").Value = 1
$ws2.Range("G34").Formula = "=F34*E34"
$ws2.Range("C35").Value = 1
$ws2.Range("E35").Formula = "=D35*C35"
$ws2.Range("F35").Value = 1
$ws2.Range("G35").Formula = "=F35*E35"
$ws2.Range("C36").Value = 1
$ws2.Range("E36").Formula = "=D36*C36"
$ws2.Range("F36").Value = 1
$ws2.Range("G36").Formula = "=F36*E36"
$ws2.Range("C37").Value = 1
$ws2.Range("E37").Formula = "=D37*C37"
$ws2.Range("F37").Value = 1
$ws2.Range("G37").Formula = "=F37*E37"
$ws2.Range("C38").Value = 1
$ws2.Range("E38").Formula = "=D38*C38"
$ws2.Range("F38").Value = 1
$ws2.Range("G38").Formula = "=F38*E38"
$ws2.Range("C39").Value = 1
$ws2.Range("E39").Formula = "=D39*C39"
$ws2.Range("F39").Value = 1
$ws2.Range("G39").Formula = "=F39*E39"
$ws2.Range("C40").Value = 1
$ws2.Range("E40").Formula = "=D40*C40"
$ws2.Range("F40").Value = 1
$ws2.Range("G40").Formula = "=F40*E40"
$ws2.Range("C41").Value = 1
$ws2.Range("E41").Formula = "=D41*C41"
$ws2.Range("F41").Value = 1
$ws2.Range("G41").Formula = "=F41*E41"
$ws2.Range("C42").Value = 1
$ws2.Range("E42").Formula = "=D42*C42"
$ws2.Range("F42").Value = 1
$ws2.Range("G42").Formula = "=F42*E42"
$ws2.Range("C43").Value = 1
$ws2.Range("E43").Formula = "=D43*C43"
$ws2.Range("F43").Value = 1
$ws2.Range("G43").Formula = "=F43*E43"
$ws2.Range("C44").Value = 1
$ws2.Range("E44").Formula = "=D44*C44"
$ws2.Range("F44").Value = 1
$ws2.Range("G44").Formula = "=F44*E44"
$ws2.Range("C45").Value = 1
$ws2.Range("E45").Formula = "=D45*C45"
$ws2.Range("F45").Value = 1
$ws2.Range("G45").Formula = "=F45*E45"

# --- Totals row 46 ---
$ws2.Range("E46").Formula = "=SUM(E3:E45)"
$ws2.Range("G46").Formula = "=SUM(G3:G45)"
# --- Hyperlinks for M5:M8 ---
$ws2.Hyperlinks.Add($ws2.Range("M5"), "https://www.postavrobota.cz/Vysilac-433MHz-ASK-antena-d663.htm") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("M6"), "https://www.postavrobota.cz/I2C-LCD-displej-znakovy-16x2-modry-d333.htm") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("M7"), "https://www.postavrobota.cz/Membranova-klavesnice-4x3-samolepici-d137.htm") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("M8"), "http://www.ebay.com/itm/PCF8574-PCF8574T-I2C-8-Bit-IO-GPIO-Expander-Module-for-Arduino-Raspberry-Pi-UK-/272432637606?var=&hash=item3f6e4026a6:m:mQFPRLxSj-zVxa61Qe6YVhQ") | Out-Null

# --- Column widths ---
$ws2.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws2.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws2.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
$ws2.Columns.Item(8).EntireColumn.AutoFit() | Out-Null
$ws2.Columns.Item(9).ColumnWidth = 10
$ws2.Columns.Item(10).ColumnWidth = 10

# --- Sheet selection / activation ---
$ws1.Activate()
$ws1.Range("H36").Select() | Out-Null

$ws2.Activate()
$ws2.Range("J30").Select() | Out-Null